# Update the "Förändrad" (Changed) date column (C) for all existing data
# rows from 2023-10-17 (45190) to 2023-10-19 (45192).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C393").Value = 45192

# The previously-last data row (393) becomes a "normal" interior row and
# picks up an explicit row height, matching every other data row.
$ws.Rows.Item(393).RowHeight = 15

# Append the new data row (394) for case A 44732-2023.
$ws.Range("A394").Value = "A 44732-2023"

$ws.Range("B394").Value = 45190
$ws.Range("B394").NumberFormat = "YYYY-MM-DD"

$ws.Range("C394").Value = 45192
$ws.Range("C394").NumberFormat = "YYYY-MM-DD"

$ws.Range("D394").Value = "SKÅNE LÄN"
$ws.Range("E394").Value = "ÖSTRA GÖINGE"

$ws.Range("G394").Value = 4.1
$ws.Range("H394:Q394").Value = 0

$ws.Range("R394").Value = ""
$ws.Range("R394").WrapText = $true

Write-Host "done"
